$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. "Evette " / "Mestre" (spell-checked) / " Meyer" -> single plain run
#    "Evette Mestre Meyer" (drops the proofErr spell-check wrapper runs).
# ---------------------------------------------------------------------------
$p1 = $d.Paragraphs.Item(1)
$nameRange = $p1.Range
$nameRange.End = $nameRange.End - 1          # exclude the paragraph mark
$nameRange.Delete()
$nameRange.InsertAfter("Evette Mestre Meyer")

# ---------------------------------------------------------------------------
# Helper: locate the Nth paragraph (1-based occurrence) whose trimmed text
# equals the given string, starting the search at paragraph index
# $startAt (1-based, inclusive). Returns the paragraph index, or 0.
# ---------------------------------------------------------------------------
function Find-ParagraphIndex($doc, [string]$text, [int]$startAt) {
    $count = $doc.Paragraphs.Count
    for ($i = $startAt; $i -le $count; $i++) {
        $para = $doc.Paragraphs.Item($i)
        $t = $para.Range.Text
        $t = $t.TrimEnd([char]13, [char]7)
        if ($t -eq $text) {
            return $i
        }
    }
    return 0
}

# ---------------------------------------------------------------------------
# 2. "Will each solution work for ALL cases?" (first occurrence, under the
#    Cat/Parrot/Seed problem) gains a trailing " Yes." answer, with "Yes."
#    colored 3366FF.
# ---------------------------------------------------------------------------
$idx = Find-ParagraphIndex $d "Will each solution work for ALL cases?" 1
$p = $d.Paragraphs.Item($idx)
$r = $p.Range
$r.End = $r.End - 1
$r.Collapse(0)
$r.InsertAfter(" ")
$r.Collapse(0)
$r.InsertAfter("Yes.")
$r.Font.Color = 16737843   # wdColor BGR encoding of RGB 3366FF

# ---------------------------------------------------------------------------
# 3. "Explain the solution in full." (first occurrence) gains a trailing
#    explanation, colored 0000FF.
# ---------------------------------------------------------------------------
$idx = Find-ParagraphIndex $d "Explain the solution in full." 1
$p = $d.Paragraphs.Item($idx)
$r = $p.Range
$r.End = $r.End - 1
$r.Collapse(0)
$r.InsertAfter(" ")
$r.Collapse(0)
$r.InsertAfter("First, secure the goat with something he can" + [char]8217 + "t chew through like a chain.  Hide the seed from the parrot.  Transport the cat across the river first.  Return for the parrot.  Then collect the cabbage and transport that across the river.  Finally, put the seed in the boat, release the goat and make your way across the river with the seed.")
$r.Font.Color = 16711680   # wdColor BGR encoding of RGB 0000FF

# ---------------------------------------------------------------------------
# 4. "Describe some test cases you tried out to make sure it works." (first
#    occurrence) gains a trailing blank run (just a space character).
# ---------------------------------------------------------------------------
$idx = Find-ParagraphIndex $d "Describe some test cases you tried out to make sure it works." 1
$p = $d.Paragraphs.Item($idx)
$r = $p.Range
$r.End = $r.End - 1
$r.Collapse(0)
$r.InsertAfter(" ")

# ---------------------------------------------------------------------------
# 5. Move the (hidden) "_GoBack" bookmark from the end of "Does each
#    solution meet the goals?" paragraph to the very start of the "Socks in
#    the Dark" heading paragraph.
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $bm = $d.Bookmarks.Item("_GoBack")
    $bm.Delete()
}
$sidx = Find-ParagraphIndex $d "Socks in the Dark" 1
$sp = $d.Paragraphs.Item($sidx)
$sr = $sp.Range
$sr.Collapse(1)   # wdCollapseStart
$d.Bookmarks.Add("_GoBack", $sr)
